$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 227.84616
$ws.Range("I9").Value = 242.22223
$ws.Range("J9").Value = 195.5
$ws.Range("K9").Value = 242.22223
$ws.Range("L9").Value = 195.5
$ws.Range("M9").Value = -73.22223
$ws.Range("N9").Value = -533.5
$ws.Range("H64").Value = 7199.8
$ws.Range("J64").Value = 8999.5
$ws.Range("L64").Value = 8999.5
$ws.Range("N64").Value = -9495.5
$ws.Range("H67").Value = 7199.8
$ws.Range("J67").Value = 8999.5
$ws.Range("L67").Value = 8999.5
$ws.Range("N67").Value = -10715.5
$ws.Range("H74").Value = 128420.664
$ws.Range("I74").Value = 143723.25
$ws.Range("K74").Value = 143723.25
$ws.Range("M74").Value = -142787.25
$ws.Range("H77").Value = 128420.664
$ws.Range("I77").Value = 143723.25
$ws.Range("K77").Value = 718616.25
$ws.Range("M77").Value = -713936.25
$ws.Range("H107").Value = 2085.2222
$ws.Range("I107").Value = 1497.75
$ws.Range("J107").Value = 2555.2
$ws.Range("K107").Value = 1497.75
$ws.Range("L107").Value = 2555.2
$ws.Range("M107").Value = 422.25
$ws.Range("N107").Value = -6395.2
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("H125").Value = 1142.5
$ws.Range("I125").Value = 999
$ws.Range("J125").Value = 1190.3334
$ws.Range("K125").Value = 8991
$ws.Range("L125").Value = 10713.0006
$ws.Range("M125").Value = -6531
$ws.Range("N125").Value = -15633.0006
$ws.Range("H131").Value = 4575544
$ws.Range("I131").Value = 33483
$ws.Range("K131").Value = 100449
$ws.Range("M131").Value = -95409
$ws.Range("H132").Value = 5550.9165
$ws.Range("I132").Value = 2555.2666
$ws.Range("J132").Value = 10543.667
$ws.Range("K132").Value = 7665.7998
$ws.Range("L132").Value = 31631.001
$ws.Range("M132").Value = -5135.7998
$ws.Range("N132").Value = -36691.001
$ws.Range("H137").Value = 2601702.8
$ws.Range("I137").Value = 4999.5
$ws.Range("J137").Value = 4332838.5
$ws.Range("K137").Value = 14998.5
$ws.Range("L137").Value = 12998515.5
$ws.Range("M137").Value = -12448.5
$ws.Range("N137").Value = -13003615.5
$ws.Range("H141").Value = 15159170
$ws.Range("I141").Value = 20836730
$ws.Range("K141").Value = 62510190
$ws.Range("M141").Value = -62505010

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2118.4314
$ws.Range("I32").Value = 2060.8
$ws.Range("K32").Value = 2060.8
$ws.Range("M32").Value = -1773.8
$ws.Range("H45").Value = 3171.75
$ws.Range("I45").Value = 1414.375
$ws.Range("J45").Value = 6686.5
$ws.Range("K45").Value = 1414.375
$ws.Range("L45").Value = 6686.5
$ws.Range("M45").Value = -1037.375
$ws.Range("N45").Value = -7440.5
$ws.Range("H61").Value = 15002297
$ws.Range("I61").Value = 22224846
$ws.Range("J61").Value = 2001709.6
$ws.Range("K61").Value = 22224846
$ws.Range("L61").Value = 2001709.6
$ws.Range("M61").Value = -22224634
$ws.Range("N61").Value = -2002133.6
$ws.Range("H62").Value = 114999.336
$ws.Range("J62").Value = 114999.336
$ws.Range("L62").Value = 114999.336
$ws.Range("N62").Value = -116247.336
$ws.Range("H65").Value = 114999.336
$ws.Range("J65").Value = 114999.336
$ws.Range("L65").Value = 344998.008
$ws.Range("N65").Value = -351238.008
$ws.Range("H74").Value = 864133.8
$ws.Range("J74").Value = 4918
$ws.Range("L74").Value = 4918
$ws.Range("N74").Value = -6666
$ws.Range("H77").Value = 864133.8
$ws.Range("J77").Value = 4918
$ws.Range("L77").Value = 24590
$ws.Range("N77").Value = -33326
$ws.Range("H88").Value = 4248.5
$ws.Range("J88").Value = 4248.5
$ws.Range("L88").Value = 4248.5
$ws.Range("N88").Value = -5060.5
$ws.Range("H91").Value = 4248.5
$ws.Range("J91").Value = 4248.5
$ws.Range("L91").Value = 4248.5
$ws.Range("N91").Value = -7056.5
$ws.Range("H94").Value = 81061.75
$ws.Range("J94").Value = 81061.75
$ws.Range("L94").Value = 81061.75
$ws.Range("N94").Value = -82863.75
$ws.Range("H122").Value = 6449.75
$ws.Range("I122").Value = 6933
$ws.Range("K122").Value = 20799
$ws.Range("M122").Value = -18349
$ws.Range("H132").Value = 3231018.8
$ws.Range("I132").Value = 5564.6523
$ws.Range("J132").Value = 12504199
$ws.Range("K132").Value = 16693.9569
$ws.Range("L132").Value = 37512597
$ws.Range("M132").Value = -14163.9569
$ws.Range("N132").Value = -37517657
$ws.Range("H136").Value = 15002297
$ws.Range("I136").Value = 22224846
$ws.Range("J136").Value = 2001709.6
$ws.Range("K136").Value = 66674538
$ws.Range("L136").Value = 6005128.800000001
$ws.Range("M136").Value = -66671988
$ws.Range("N136").Value = -6010228.800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2112.325
$ws.Range("I94").Value = 1716.4828
$ws.Range("J94").Value = 3155.9092
$ws.Range("K94").Value = 1716.4828
$ws.Range("L94").Value = 3155.9092
$ws.Range("M94").Value = -1265.4828
$ws.Range("N94").Value = -4057.9092
$ws.Range("H134").Value = 7145846
$ws.Range("I134").Value = 2872.2222
$ws.Range("J134").Value = 20003198
$ws.Range("K134").Value = 8616.6666
$ws.Range("L134").Value = 60009594
$ws.Range("M134").Value = -6081.6666
$ws.Range("N134").Value = -60014664

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36082330
$ws.Range("I31").Value = 41670696
$ws.Range("K31").Value = 41670696
$ws.Range("M31").Value = -41670401
$ws.Range("H34").Value = 36082330
$ws.Range("I34").Value = 41670696
$ws.Range("K34").Value = 41670696
$ws.Range("M34").Value = -41670494
$ws.Range("H58").Value = 2382.5
$ws.Range("I58").Value = 2274.6667
$ws.Range("K58").Value = 2274.6667
$ws.Range("M58").Value = -2071.6667
$ws.Range("H92").Value = 31624.75
$ws.Range("J92").Value = 31624.75
$ws.Range("L92").Value = 31624.75
$ws.Range("N92").Value = -36616.75
$ws.Range("H99").Value = 24119.824
$ws.Range("I99").Value = 7928.75
$ws.Range("K99").Value = 7928.75
$ws.Range("M99").Value = -6430.75
$ws.Range("H107").Value = 2562.8918
$ws.Range("I107").Value = 2470.88
$ws.Range("J107").Value = 2754.5833
$ws.Range("K107").Value = 2470.88
$ws.Range("L107").Value = 2754.5833
$ws.Range("M107").Value = -550.8800000000001
$ws.Range("N107").Value = -6594.5833
$ws.Range("H122").Value = 3501.2
$ws.Range("I122").Value = 3314
$ws.Range("K122").Value = 9942
$ws.Range("M122").Value = -7492
$ws.Range("H126").Value = 24119.824
$ws.Range("I126").Value = 7928.75
$ws.Range("K126").Value = 23786.25
$ws.Range("M126").Value = -21316.25
$ws.Range("H132").Value = 3144.375
$ws.Range("I132").Value = 2859.4167
$ws.Range("J132").Value = 3999.25
$ws.Range("K132").Value = 8578.250100000001
$ws.Range("L132").Value = 11997.75
$ws.Range("M132").Value = -6048.250100000001
$ws.Range("N132").Value = -17057.75
$ws.Range("H134").Value = 2057.2354
$ws.Range("I134").Value = 2130.3845
$ws.Range("J134").Value = 1819.5
$ws.Range("K134").Value = 6391.1535
$ws.Range("L134").Value = 5458.5
$ws.Range("M134").Value = -3856.1535
$ws.Range("N134").Value = -10528.5
$ws.Range("H136").Value = 2382.5
$ws.Range("I136").Value = 2274.6667
$ws.Range("K136").Value = 6824.000100000001
$ws.Range("M136").Value = -4274.000100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6244.615
$ws.Range("J131").Value = 7018.75
$ws.Range("L131").Value = 21056.25
$ws.Range("N131").Value = -31136.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11233527
$ws.Range("I132").Value = 3938.7273
$ws.Range("J132").Value = 31821106
$ws.Range("K132").Value = 11816.1819
$ws.Range("L132").Value = 95463318
$ws.Range("M132").Value = -9286.1819
$ws.Range("N132").Value = -95468378

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4171.6665
$ws.Range("I16").Value = 2007.1818
$ws.Range("J16").Value = 10124
$ws.Range("K16").Value = 2007.1818
$ws.Range("L16").Value = 10124
$ws.Range("M16").Value = -1837.1818
$ws.Range("N16").Value = -10464
$ws.Range("H22").Value = 4254.4443
$ws.Range("J22").Value = 4398.4
$ws.Range("L22").Value = 4398.4
$ws.Range("N22").Value = -4988.4
$ws.Range("H27").Value = 4254.4443
$ws.Range("J27").Value = 4398.4
$ws.Range("L27").Value = 4398.4
$ws.Range("N27").Value = -4612.4
$ws.Range("H40").Value = 6098.7144
$ws.Range("I40").Value = 6031.8335
$ws.Range("K40").Value = 6031.8335
$ws.Range("M40").Value = -5895.8335
$ws.Range("H55").Value = 1448.579
$ws.Range("I55").Value = 2069.8
$ws.Range("J55").Value = 1226.7142
$ws.Range("K55").Value = 2069.8
$ws.Range("L55").Value = 1226.7142
$ws.Range("M55").Value = -1896.8
$ws.Range("N55").Value = -1572.7142
$ws.Range("H56").Value = 50012.5
$ws.Range("I56").Value = 50025.5
$ws.Range("J56").Value = 49999.5
$ws.Range("K56").Value = 50025.5
$ws.Range("L56").Value = 49999.5
$ws.Range("M56").Value = -49334.5
$ws.Range("N56").Value = -51381.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("H136").Value = 3054.923
$ws.Range("I136").Value = 1973.7778
$ws.Range("K136").Value = 5921.3334
$ws.Range("M136").Value = -3371.3334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 418414.03
$ws.Range("I132").Value = 1524.3889
$ws.Range("J132").Value = 1669083
$ws.Range("K132").Value = 4573.1667
$ws.Range("L132").Value = 5007249
$ws.Range("M132").Value = -2043.1667
$ws.Range("N132").Value = -5012309

# ---- Special cases: row structure changed, trailing cell becomes fully empty ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N121").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N118").Value = ""
